$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> { Column -> new value } updates, derived from the source diff.
# Each new value is written with a leading apostrophe (forces literal text,
# matching the original cells which are stored as text/inline strings) and
# the cells Style is restored afterward so no spurious "Text" number format
# gets attached (keeps the style index identical to the original workbook).
$updates = @{
    2 = @{ D='27.038.85'; E='  +0.65%  ' }
    3 = @{ D='1.680.02'; E='  +0.81%  ' }
    4 = @{ E='  +0.04%  ' }
    5 = @{ D='215.97'; E='  +0.27%  ' }
    6 = @{ E='  -3.15%  ' }
    8 = @{ D='21.67'; E='  +7.38%  ' }
    9 = @{ E='  +0.92%  ' }
    10 = @{ E='  +0.50%  ' }
    11 = @{ E='  -0.67%  ' }
    12 = @{ D='1.916.82'; E='  +0.78%  ' }
    13 = @{ D='1.671.26'; E='  +0.84%  ' }
    14 = @{ E='  +0.48%  ' }
    15 = @{ E='  +1.61%  ' }
    16 = @{ D='66.51'; E='  +0.80%  ' }
    17 = @{ D='27.026.38'; E='  +0.58%  ' }
    18 = @{ D='8.19'; E='  +5.15%  ' }
    19 = @{ D='235.77'; E='  +1.65%  ' }
    20 = @{ D='0.0₃0740'; E='  +1.02%  ' }
    22 = @{ D='4.47'; E='  +0.14%  ' }
    23 = @{ E='  +0.97%  ' }
    24 = @{ E='  -4.11%  ' }
    25 = @{ D='146.59' }
    26 = @{ E='  +1.75%  ' }
    27 = @{ D='16.47'; E='  +3.68%  ' }
    28 = @{ E='  -2.69%  ' }
    29 = @{ E='  +0.03%  ' }
    30 = @{ E='  +0.59%  ' }
    31 = @{ E='  -0.16%  ' }
    32 = @{ E='  +0.11%  ' }
    33 = @{ D='1.520.84'; E='  +3.56%  ' }
    34 = @{ E='  +0.43%  ' }
    35 = @{ D='1.71'; E='  +4.85%  ' }
    36 = @{ E='  -0.67%  ' }
    37 = @{ D='0.593'; E='  +3.36%  ' }
    38 = @{ D='0.921'; E='  +2.73%  ' }
    39 = @{ E='  +3.21%  ' }
    40 = @{ E='  +5.68%  ' }
    41 = @{ D='5.71'; E='  -2.54%  ' }
    43 = @{ D='67.90'; E='  +3.05%  ' }
    44 = @{ E='  -0.32%  ' }
    45 = @{ D='1.822.08'; E='  +0.42%  ' }
    46 = @{ E='  +0.54%  ' }
    47 = @{ D='90.21'; E='  -0.07%  ' }
    48 = @{ E='  -0.24%  ' }
    49 = @{ E='  +2.79%  ' }
    50 = @{ E='  +5.42%  ' }
    51 = @{ E='  -0.46%  ' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $origStyle = $cell.Style
        $cell.Value = "'" + $cols[$col]
        $cell.Style = $origStyle
    }
}
